$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh price / 1h-volume figures (and the two swapped coin rows) to match the latest pull.
$ws.Range("D2").Value = '''42.919.88'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '''2.529.32'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''317.27'
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").Value = '''95.64'
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.533'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '''36.09'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '''7.59'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = '''2.922.24'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '''2.526.04'
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").Value = '''15.33'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '''0.849'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '''43.012.87'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '''12.98'
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = '''6.65'
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").Value = '''0.0₃0966'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("D22").Value = '''70.21'
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = '''251.70'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").Value = '''26.90'
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '''2.41'
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("D29").Value = '''39.91'
$ws.Range("E29").Value = '  +4.39%  '
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").Value = '''154.86'
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("E33").Value = '  +2.70%  '
$ws.Range("D34").Value = '''3.30'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").Value = '''0.0790'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").Value = '''18.85'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = '''0.112'
$ws.Range("E38").Value = '  -3.05%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '''23.77'
$ws.Range("E40").Value = '  -1.23%  '
$ws.Range("D41").Value = '''2.30'
$ws.Range("E41").Value = '  +11.66%  '
$ws.Range("D42").Value = '''0.0304'
$ws.Range("E42").Value = '  +1.24%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''3.79'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("E45").Value = '  -3.08%  '
$ws.Range("D46").Value = '''2.020.73'
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("D49").Value = '''2.776.28'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").Value = '''73.62'
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").Value = '''102.57'
$ws.Range("E51").Value = '  +0.78%  '
